{"js": "// Replace the whole body content with the updated Python test-case script.\n// The document is rewritten from a password-validation snippet (with a manual\n// test harness / re-based validation) into a proper unittest-based test\n// suite plus a simplified validate_password() implementation, wrapped in a\n// ```python code fence.\nconst body = context.document.body;\n\n// Start from a clean slate.\nbody.clear();\nawait context.sync();\n\nconst newText = `\\`\\`\\`python\nimport unittest\n\nclass TestPasswordValidation(unittest.TestCase):\n\n    def test_password_too_short(self):\n        result = validate_password(\"abc123\")\n        self.assertEqual(result, \"Password must be at least 8 characters.\")\n\n    def test_password_no_number(self):\n        result = validate_password(\"abcdefgh\")\n        self.assertEqual(result, \"Password must include at least one number.\")\n\n    def test_password_short_with_number(self):\n        result = validate_password(\"abc1234\")\n        self.assertEqual(result, \"Password must be at least 8 characters.\")\n\n    def test_valid_password(self):\n        result = validate_password(\"abc12345\")\n        self.assertEqual(result, True)\n\n    def test_long_password_with_numbers(self):\n        result = validate_password(\"mypassword1\")\n        self.assertEqual(result, True)\n\n    def test_all_numbers_valid(self):\n        result = validate_password(\"12345678\")\n        self.assertEqual(result, True)\n\ndef validate_password(password):\n    if len(password) < 8:\n        return \"Password must be at least 8 characters.\"\n    if not any(char.isdigit() for char in password):\n        return \"Password must include at least one number.\"\n    return True\n\nif __name__ == '__main__':\n    unittest.main()\n\\`\\`\\``;\n\n// insertText() with \"\\n\"-separated lines creates one paragraph per line, so\n// blank lines in the template literal become the empty paragraphs that\n// separate the code blocks below, matching the target layout exactly.\nbody.insertText(newText, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Replace the whole document body with the updated Python test-case script:\n# rewrite the password-validation snippet (and its manual test harness) into\n# a proper unittest-based test suite plus a simplified validate_password()\n# implementation, wrapped in a ```python code fence.\n$d = $word.ActiveDocument\n\n$newBody = @'\n```python\nimport unittest\n\nclass TestPasswordValidation(unittest.TestCase):\n\n    def test_password_too_short(self):\n        result = validate_password(\"abc123\")\n        self.assertEqual(result, \"Password must be at least 8 characters.\")\n\n    def test_password_no_number(self):\n        result = validate_password(\"abcdefgh\")\n        self.assertEqual(result, \"Password must include at least one number.\")\n\n    def test_password_short_with_number(self):\n        result = validate_password(\"abc1234\")\n        self.assertEqual(result, \"Password must be at least 8 characters.\")\n\n    def test_valid_password(self):\n        result = validate_password(\"abc12345\")\n        self.assertEqual(result, True)\n\n    def test_long_password_with_numbers(self):\n        result = validate_password(\"mypassword1\")\n        self.assertEqual(result, True)\n\n    def test_all_numbers_valid(self):\n        result = validate_password(\"12345678\")\n        self.assertEqual(result, True)\n\ndef validate_password(password):\n    if len(password) < 8:\n        return \"Password must be at least 8 characters.\"\n    if not any(char.isdigit() for char in password):\n        return \"Password must include at least one number.\"\n    return True\n\nif __name__ == '__main__':\n    unittest.main()\n```\n'@\n\n# Here-strings only use \"`n\" between lines; Word paragraph marks are \"`r\",\n# so convert before assigning.\n$newBody = $newBody -replace \"`n\", \"`r\"\n\n$full = $d.Range(0, $d.Content.End)\n$full.Text = $newBody\n"}
